$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 810 ("「冬で私が好きな物」" post), shifting all rows below it up by one.
$ws.Rows.Item(810).Delete()
